# feature/0020 - Release 1.3
# Add a new "1.3" worksheet (clone of the "1.2" QA checklist) at the end of
# the workbook, carrying over the same values/formatting but without the
# custom column width that "1.2" had picked up. Then update the view state
# (active tab / selections) on both sheets to match what Excel leaves behind
# after such a sheet-duplication + tab switch.

$wb = $excel.ActiveWorkbook

# "1.2" is the 4th tab and is the template for the new "1.3" tab.
$src = $wb.Worksheets.Item(4)

# Insert the new sheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$new = $wb.Worksheets.Add($null, $lastSheet)
$new.Name = "1.3"

# --- Header row (B1:D1) ---
$new.Range("B1").Value = $src.Range("B1").Value2
$new.Range("C1").Value = $src.Range("C1").Value2
$new.Range("D1").Value = $src.Range("D1").Value2

# --- Data rows (A2:A23 labels, B:D formatted "Good"-style answer cells) ---
$labels = @(
    "Search (modal)",
    "Product list",
    "Product list recommend",
    "Product detail",
    "HP banners",
    "HP carousels",
    "Basket (float)",
    "Checkout basket",
    "Checkout recommend",
    "Favourite",
    "Skrývá se náhodný produkt?",
    "Neovlivňuje náhodý produkt jiné?",
    "Skrývá se náhodný produkt ze specific?",
    "Neovlivňuje náhodý produkt jiné ze specific?",
    "Zobrazuje se produkt s false?",
    "Zobrazuje se produkt ze specific s false?",
    "Funguje default action",
    "Funguje hide action",
    "Funguje mark action",
    "Funguje transparent action",
    "Lze změnit action",
    "Funguje neverskip?"
)

for ($i = 0; $i -lt $labels.Count; $i++) {
    $row = $i + 2
    $new.Range("A$row").Value = $labels[$i]
    $rowRange = $new.Range("B$row`:D$row")
    $rowRange.Style = "Good"
}

# --- View state ---
# "1.2" keeps its data but now shows the whole table selected (no more
# single active-cell selection) and scrolled down near the bottom.
$src.Activate()
$src.Range("A1:D23").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1

# The new "1.3" tab becomes the active / displayed tab, with the whole
# table selected too.
$new.Activate()
$new.Range("A1:D23").Select()
